$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    "K" = 1
    "L" = 5
    "O" = 5.7680816650390625
    "P" = 2.2352011203765869
    "Q" = 2.2352011203765869
    "R" = 0
    "S" = 19.658470153808594
    "T" = 42.728260040283203
    "U" = 40.030673980712891
    "V" = 2.6975860595703125
    "X" = 1
    "Y" = 16.939168930053711
    "Z" = 16.453506469726563
    "AA" = 16.453506469726563
    "AB" = 0
    "AC" = 14.110039710998535
    "AD" = 34.581287384033203
    "AE" = 7.0867562294006348
    "AF" = 27.494531631469727
    "AG" = 0
    "AI" = 4.0451245307922363
    "AJ" = 17.451549530029297
    "AK" = 17.451549530029297
    "AM" = 2.4514777660369873
    "AN" = 44.512004852294922
    "AO" = 17.808340072631836
    "AP" = 26.703664779663086
    "AQ" = 0
    "AS" = 3.6109740734100342
    "AT" = 27.876899719238281
    "AU" = 20.809535980224609
    "AV" = 7.0673637390136719
    "AW" = 9.8803443908691406
    "AX" = 8.5362758636474609
    "AY" = 8.5362758636474609
    "BA" = 0
    "BC" = 3.2463028430938721
    "BD" = 32.966930389404297
    "BE" = 32.966930389404297
    "BF" = 0
    "BG" = 8.0707817077636719
    "BH" = 11.795376777648926
    "BI" = 11.795376777648926
    "BJ" = 0
    "BM" = 5.4777421951293945
    "BN" = 44.768749237060547
    "BO" = 5.536811351776123
    "BP" = 39.231937408447266
    "BQ" = 1.8907539844512939
    "BR" = 6.6376943588256836
    "BS" = 3.558706521987915
    "BT" = 3.0789878368377686
    "BU" = 10.387241363525391
    "BV" = 7.8235669136047363
}

foreach ($col in $newValues.Keys) {
    $value = $newValues[$col]
    $ws.Range("${col}2").Value = $value
    $ws.Range("${col}3").Value = $value
}
